$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like text in column D stays text, not auto-converted to a date serial
$ws.Range('D919:D948').NumberFormat = '@'

# Row 919
$ws.Cells.Item(919, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(919, 2).Value = 'Back to School, Back to Basics'
$ws.Cells.Item(919, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2019/08/back_to_school_back.html'
$ws.Cells.Item(919, 4).Value = '2023-02-14'
$ws.Cells.Item(919, 5).Value = 'More than 300 million students yearly attend schools, colleges and universities worldwide.'

# Row 920
$ws.Cells.Item(920, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(920, 2).Value = 'Roses and Chocolates and Bears, Oh My!'
$ws.Cells.Item(920, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2019/02/roses_and_chocolates.html'
$ws.Cells.Item(920, 4).Value = '2023-02-14'
$ws.Cells.Item(920, 5).Value = 'From flowers, to confectionary products, to stuffed animals and more, consumers spent millions of dollars worldwide to celebrate Valentine’s Day, making it one of the largest retail shopping days of the year.'

# Row 921
$ws.Cells.Item(921, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(921, 2).Value = 'Annotating an Export Shipment: Filing Citations, Exemption and Exclusion Legends'
$ws.Cells.Item(921, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2016/10/annotating_an_export.html'
$ws.Cells.Item(921, 4).Value = '2021-10-08'
$ws.Cells.Item(921, 5).Value = 'The U.S. Census Bureau often receives questions on how to annotate commercial documents for export shipments to minimize potential delays at the port of export.'

# Row 922
$ws.Cells.Item(922, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(922, 2).Value = 'February 2015, Trade Deficit Decreased'
$ws.Cells.Item(922, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2015/04/february-2015-trade-deficit-decreased.html'
$ws.Cells.Item(922, 4).Value = '2021-12-16'
$ws.Cells.Item(922, 5).Value = 'The trade deficit in goods and services decreased to $35.4 billion in February, a $7.2 billion decrease from the January deficit of $42.7 billion, revised.'

# Row 923
$ws.Cells.Item(923, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(923, 2).Value = 'January 2015, Trade Deficit Decreased'
$ws.Cells.Item(923, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2015/03/january-2015-trade-deficit-decreased.html'
$ws.Cells.Item(923, 4).Value = '2021-12-16'
$ws.Cells.Item(923, 5).Value = 'The trade deficit in goods and services decreased to $41.8 billion in January, a $3.8 billion decrease from the December deficit of $45.6 billion, revised.'

# Row 924
$ws.Cells.Item(924, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(924, 2).Value = 'December 2014, Trade Deficit Increased'
$ws.Cells.Item(924, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2015/02/december-2014-trade-deficit-increased.html'
$ws.Cells.Item(924, 4).Value = '2021-12-16'
$ws.Cells.Item(924, 5).Value = 'The trade deficit in goods and services increased by the highest margin on record ($6.8 billion) to $46.6 billion in December, a 17.1% increase from November ($39.8 billion, revised).'

# Row 925
$ws.Cells.Item(925, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(925, 2).Value = 'November 2014, Trade Deficit Decreased'
$ws.Cells.Item(925, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2015/01/november-2014-trade-deficit-decreased.html'
$ws.Cells.Item(925, 4).Value = '2021-12-16'
$ws.Cells.Item(925, 5).Value = 'The trade deficit in goods and services continued to decrease. In November, the deficit was $39.0 billion, down $3.2 billion from October.'

# Row 926
$ws.Cells.Item(926, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(926, 2).Value = 'October 2014, Trade Deficit Decrease'
$ws.Cells.Item(926, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2014/12/october-2014-trade-deficit-decrease.html'
$ws.Cells.Item(926, 4).Value = '2021-12-16'
$ws.Cells.Item(926, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 927
$ws.Cells.Item(927, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(927, 2).Value = 'September 2014, Trade Deficit Increase'
$ws.Cells.Item(927, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2014/11/september-2014-trade-deficit-increase.html'
$ws.Cells.Item(927, 4).Value = '2021-12-16'
$ws.Cells.Item(927, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 928
$ws.Cells.Item(928, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(928, 2).Value = 'August 2014, Trade Deficit Decreases'
$ws.Cells.Item(928, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2014/10/august-2014-trade-deficit-decreases.html'
$ws.Cells.Item(928, 4).Value = '2021-10-08'
$ws.Cells.Item(928, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 929
$ws.Cells.Item(929, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(929, 2).Value = 'July 2014, Trade Deficit Decreases'
$ws.Cells.Item(929, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2014/09/july-2014-trade-deficit-decreases.html'
$ws.Cells.Item(929, 4).Value = '2021-12-16'
$ws.Cells.Item(929, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 930
$ws.Cells.Item(930, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(930, 2).Value = 'June 2014, Second Month in a Row Trade Deficit Decreases'
$ws.Cells.Item(930, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2014/08/june-2014-second-month-in-a-row-trade-deficit-decreases.html'
$ws.Cells.Item(930, 4).Value = '2021-12-16'
$ws.Cells.Item(930, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 931
$ws.Cells.Item(931, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(931, 2).Value = 'New Year, New Export Markets!'
$ws.Cells.Item(931, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2014/01/new-year-new-export-markets.html'
$ws.Cells.Item(931, 4).Value = '2021-12-16'
$ws.Cells.Item(931, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 932
$ws.Cells.Item(932, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(932, 2).Value = 'Find Buyers, Make Contacts Abroad with Foreign Agricultural Service'
$ws.Cells.Item(932, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2013/07/fas_webinar.html'
$ws.Cells.Item(932, 4).Value = '2021-12-16'
$ws.Cells.Item(932, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 933
$ws.Cells.Item(933, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(933, 2).Value = 'Simple Question, Big Impact � What is Your Port of Export?'
$ws.Cells.Item(933, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2013/03/simple-question-big-impact-what-is-your-port-of-export.html'
$ws.Cells.Item(933, 4).Value = '2021-12-16'
$ws.Cells.Item(933, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 934
$ws.Cells.Item(934, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(934, 2).Value = 'Deficit Balance Increased, Exports and Imports Decreased in Oct. 2012'
$ws.Cells.Item(934, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2012/12/october-2012-trade.html'
$ws.Cells.Item(934, 4).Value = '2022-05-17'
$ws.Cells.Item(934, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 935
$ws.Cells.Item(935, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(935, 2).Value = 'Survey of Business Owners, Meet the Profile.'
$ws.Cells.Item(935, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2012/10/survey-of-business-profile.html'
$ws.Cells.Item(935, 4).Value = '2021-12-16'
$ws.Cells.Item(935, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 936
$ws.Cells.Item(936, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(936, 2).Value = 'Back to the Basics – Ensuring Your Privacy, Leave the EIN Out!'
$ws.Cells.Item(936, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2012/05/back-to-the-basics.html'
$ws.Cells.Item(936, 4).Value = '2022-09-09'
$ws.Cells.Item(936, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 937
$ws.Cells.Item(937, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(937, 2).Value = 'globalEDGE: A World of Resources, A World of Opportunities!'
$ws.Cells.Item(937, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2012/02/globaledge-a-world-of-resources-a-world-of-opportunities.html'
$ws.Cells.Item(937, 4).Value = '2021-12-16'
$ws.Cells.Item(937, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 938
$ws.Cells.Item(938, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(938, 2).Value = 'AESDirect is Moving to a New, Improved Platform!'
$ws.Cells.Item(938, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2011/08/aesdirect-is-moving-to-a-new-improved-platform.html'
$ws.Cells.Item(938, 4).Value = '2022-04-21'
$ws.Cells.Item(938, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 939
$ws.Cells.Item(939, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(939, 2).Value = 'Exporting Software: To file or not to file, Part 1'
$ws.Cells.Item(939, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2010/11/exporting-software-to-file-or-not-to-file-part-1.html'
$ws.Cells.Item(939, 4).Value = '2021-12-16'
$ws.Cells.Item(939, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 940
$ws.Cells.Item(940, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(940, 2).Value = 'Where Are Your Goods From, Originally?'
$ws.Cells.Item(940, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2010/10/where-are-your-goods-from-originally.html'
$ws.Cells.Item(940, 4).Value = '2021-12-16'
$ws.Cells.Item(940, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 941
$ws.Cells.Item(941, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(941, 2).Value = 'Real Export Emergencies, Episode 3: Send an SOS to AES'
$ws.Cells.Item(941, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2010/08/real-export-emergencies-episode-3-send-an-sos-to-aes.html'
$ws.Cells.Item(941, 4).Value = '2021-10-08'
$ws.Cells.Item(941, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 942
$ws.Cells.Item(942, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(942, 2).Value = 'Real Export Emergencies, Episode 2: Using the Same Name for Twin Shipments'
$ws.Cells.Item(942, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2010/06/real-export-emergencies-episode-2-using-the-same-name-for-twin-shipments.html'
$ws.Cells.Item(942, 4).Value = '2021-10-08'
$ws.Cells.Item(942, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 943
$ws.Cells.Item(943, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(943, 2).Value = 'Appendices A,D,F, and U: Get to know them better'
$ws.Cells.Item(943, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2010/06/appendices-adf-and-u-get-to-know-them-better.html'
$ws.Cells.Item(943, 4).Value = '2022-04-14'
$ws.Cells.Item(943, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 944
$ws.Cells.Item(944, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(944, 2).Value = 'Why does AESDirect say: Registration Already Exists with ID Number?'
$ws.Cells.Item(944, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2010/05/why-does-aesdirect-say-registration-already-exists-with-id-number-ive-never-used-aesdirect-before-1.html'
$ws.Cells.Item(944, 4).Value = '2022-09-15'
$ws.Cells.Item(944, 5).Value = 'I tried to register through AESDirect and it says my account is already registered. What am I supposed to do?'

# Row 945
$ws.Cells.Item(945, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(945, 2).Value = 'It''s Not What You Say, It''s What You Can Prove'
$ws.Cells.Item(945, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2010/02/its-not-what-you-say-its-what-you-can-prove.html'
$ws.Cells.Item(945, 4).Value = '2022-09-12'
$ws.Cells.Item(945, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 946
$ws.Cells.Item(946, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(946, 2).Value = 'Renewable Energy, Employment and Foreign Trade'
$ws.Cells.Item(946, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2010/01/renewable-energy-employment-and-foreign-trade.html'
$ws.Cells.Item(946, 4).Value = '2021-10-08'
$ws.Cells.Item(946, 5).Value = 'Get foreign trade data reports, information on Trade Regulations and answers to your questions from the official source of US Import & Export Trade Statistics.'

# Row 947
$ws.Cells.Item(947, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(947, 2).Value = 'U.S. Exports by Metropolitan Area'
$ws.Cells.Item(947, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2018/07/u_s_exports_by_metr.html'
$ws.Cells.Item(947, 4).Value = '2021-10-08'
$ws.Cells.Item(947, 5).Value = ''

# Row 948
$ws.Cells.Item(948, 1).Value = 'Censys Global Reach'
$ws.Cells.Item(948, 2).Value = 'The Goods and Services Deficit Decreased to $43.7 billion in March 2017'
$ws.Cells.Item(948, 3).Value = 'https://www.census.gov/newsroom/blogs/global-reach/2017/05/the_goods_and_servic.html'
$ws.Cells.Item(948, 4).Value = '2021-10-08'
$ws.Cells.Item(948, 5).Value = ''
